$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C6: was a numeric value, now a text string "thisisnotanid"
$ws.Range("C6").Value = "thisisnotanid"

# Update F9: was a numeric date value, now a text string "thisisnotadate"
$ws.Range("F9").Value = "thisisnotadate"

# Update selected cell in the sheet view
$ws.Range("J10").Select()
